$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Nuevas proxys / filas en el diagrama: Monto?, FormaPago?, Faltan:
# (se insertan en este orden para que los shared-strings queden Monto?, FormaPago?, Faltan:)
$ws.Range("B12").Value = "Monto?"
$ws.Range("B13").Value = "FormaPago?"
$ws.Range("A12").Value = "Faltan:"

# Estilo en negrita para las nuevas etiquetas
$lbl = $ws.Range("B12:B13")
$lbl.Font.Bold = $true

# Correccion: Recibo pasa a color verde (clase proxy), igual que el resto de la fila 7/8
$ws.Range("D6").Font.Color = 5287936

# Ultima celda seleccionada al guardar
$ws.Range("C14").Select() | Out-Null

Write-Host "edit applied"
